$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '98.026.67'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.64%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.423.95'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.05%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '258.98'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.87%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '660.98'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.75%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.48'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -8.24%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.442'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +8.00%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.08'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.04%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '3.422.43'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.01%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.215'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.08%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '42.46'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.02%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.46'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +16.61%  '

$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '97.721.96'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.62%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000266'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +5.71%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.064.98'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.01%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '9.56'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +36.87%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.419.88'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.96%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.16'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +7.56%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.512'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +23.50%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.00'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +8.99%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.50'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.73%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '516.74'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.55%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000208'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.34%  '

$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.49'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +7.17%  '

$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '101.54'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +10.00%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '13.20'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.24%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.602.63'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.12%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.155'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.55%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '11.70'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.07%  '

$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.201'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.81%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.12%  '

$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.591'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +12.75%  '

$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.54%  '

$ws.Range('B36').Value = 'PancakeSwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.37'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +16.40%  '

$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '30.11'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.78%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '7.93'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.95%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.47'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +10.06%  '

$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '535.86'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.45%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.157'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.46%  '

$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.02%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.885'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +5.99%  '

$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '24.75'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.11%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0431'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +11.97%  '

$ws.Range('B46').Value = 'Cosmos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '8.86'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +14.12%  '

$ws.Range('B47').Value = 'MantraDAO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.70'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.46%  '

$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.41'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.41%  '

$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.61'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +12.86%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.64'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +10.18%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.10'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +4.02%  '
